$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.309.20'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.608.69'
$ws.Range("E3").Value = '  +8.01%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.89'
$ws.Range("E5").Value = '  +3.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.61'
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.601'
$ws.Range("E7").Value = '  +6.66%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").Value = '  +14.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.79'
$ws.Range("E10").Value = '  +12.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("E11").Value = '  +6.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.97'
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.32'
$ws.Range("E13").Value = '  +16.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.009.62'
$ws.Range("E14").Value = '  +8.15%  '
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.608.34'
$ws.Range("E16").Value = '  +8.30%  '
$ws.Range("E17").Value = '  +9.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.07'
$ws.Range("E18").Value = '  +6.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '46.517.02'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.35'
$ws.Range("E20").Value = '  +2.98%  '
$ws.Range("E21").Value = '  +7.22%  '
$ws.Range("E22").Value = '  +10.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.11'
$ws.Range("E23").Value = '  +5.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '255.14'
$ws.Range("E24").Value = '  +4.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.08'
$ws.Range("E25").Value = '  +10.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  +15.10%  '
$ws.Range("E27").Value = '  +32.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.60'
$ws.Range("E29").Value = '  +8.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '40.64'
$ws.Range("E30").Value = '  +3.61%  '
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.20'
$ws.Range("E32").Value = '  +11.30%  '
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("E34").Value = '  +16.27%  '
$ws.Range("E35").Value = '  +5.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0841'
$ws.Range("E36").Value = '  +8.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.07'
$ws.Range("E37").Value = '  +3.63%  '
$ws.Range("E38").Value = '  +5.55%  '
$ws.Range("E39").Value = '  +5.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.08'
$ws.Range("E40").Value = '  +10.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.25'
$ws.Range("E41").Value = '  +9.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.63'
$ws.Range("E42").Value = '  +10.92%  '
$ws.Range("E43").Value = '  +9.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.11'
$ws.Range("E44").Value = '  +45.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.035.81'
$ws.Range("E45").Value = '  +4.07%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '91.40'
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.30'
$ws.Range("E48").Value = '  +11.49%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.26'
$ws.Range("E49").Value = '  +7.05%  '
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.202'
$ws.Range("E51").Value = '  +9.05%  '
